$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed cryptocurrency price/volume data (and a Kaspa/WrappedeETH row-order
# swap in rows 27-28). Values are kept as plain text (NumberFormat "@") because
# the source sheet stores Price/Volume columns as text, and several of the new
# values (e.g. "541.40", "0.163") would otherwise be auto-coerced to numbers.
$updates = [ordered]@{
    "D2" = "59.583.11"
    "E2" = "  -4.48%  "
    "D3" = "2.489.99"
    "E3" = "  -4.89%  "
    "E4" = "  +0.12%  "
    "D5" = "541.40"
    "E5" = "  -1.92%  "
    "D6" = "147.01"
    "E6" = "  -5.12%  "
    "D7" = "0.997"
    "E7" = "  -0.22%  "
    "D8" = "0.578"
    "E8" = "  -2.67%  "
    "D9" = "2.517.58"
    "E9" = "  -3.86%  "
    "E10" = "  -3.57%  "
    "E11" = "  -1.29%  "
    "E12" = "  -0.92%  "
    "D13" = "0.358"
    "E13" = "  -1.97%  "
    "D14" = "2.933.06"
    "E14" = "  -4.73%  "
    "D15" = "24.51"
    "E15" = "  -4.64%  "
    "D16" = "59.586.30"
    "E16" = "  -4.24%  "
    "E17" = "  -2.82%  "
    "D18" = "2.512.22"
    "E18" = "  -3.95%  "
    "D19" = "11.57"
    "E19" = "  -0.75%  "
    "E20" = "  -3.42%  "
    "D21" = "326.49"
    "E21" = "  -4.29%  "
    "E22" = "  +0.37%  "
    "D23" = "5.82"
    "E23" = "  -4.69%  "
    "D24" = "61.23"
    "E24" = "  -2.81%  "
    "E25" = "  -10.05%  "
    "E26" = "  +1.02%  "
    "B27" = "Kaspa"
    "C27" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D27" = "0.163"
    "E27" = "  -2.89%  "
    "B28" = "WrappedeETH"
    "C28" = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
    "D28" = "2.619.87"
    "E28" = "  -3.82%  "
    "D29" = "7.86"
    "E29" = "  -2.05%  "
    "D30" = "7.22"
    "E30" = "  +0.73%  "
    "D31" = "0.0₃0792"
    "E31" = "  -4.35%  "
    "D32" = "1.28"
    "E32" = "  -3.99%  "
    "E33" = "  -3.69%  "
    "D34" = "159.72"
    "E34" = "  -0.57%  "
    "D35" = "0.996"
    "E35" = "  -0.32%  "
    "D36" = "1.43"
    "E36" = "  +1.04%  "
    "D37" = "18.77"
    "E37" = "  -2.60%  "
    "D38" = "4.51"
    "E38" = "  -4.03%  "
    "D39" = "1.67"
    "E39" = "  -3.73%  "
    "D40" = "6.01"
    "E40" = "  -1.73%  "
    "D41" = "314.59"
    "E41" = "  -6.58%  "
    "D42" = "36.72"
    "E42" = "  -2.70%  "
    "E43" = "  -2.86%  "
    "D44" = "0.836"
    "E44" = "  -6.23%  "
    "D45" = "0.995"
    "E45" = "  -0.37%  "
    "D46" = "0.601"
    "E46" = "  -1.67%  "
    "D47" = "10.80"
    "E47" = "  -1.50%  "
    "D48" = "125.60"
    "E48" = "  -0.51%  "
    "D49" = "0.0942"
    "E49" = "  -2.55%  "
    "D50" = "0.0529"
    "E50" = "  -3.50%  "
    "D51" = "0.0232"
    "E51" = "  -3.14%  "
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
